$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" A1 conversion summary text ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.1 = 37060.17 pesos`n✅ 37060.17 pesos = 9.07 = 951.93 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update "tasas" sheet rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 109.85
$wsTasas.Range("O10").Value = 4071.06
$wsTasas.Range("N12").Value = 4084
$wsTasas.Range("O12").Value = 104.902
